$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new columns before the existing "ExpPoints" column (C),
# pushing it to column G, and carrying over the header style.
$ws.Range("C1:F1").EntireColumn.Insert()

# New header labels for the inserted columns.
$ws.Range("C1").Value = "WIN"
$ws.Range("D1").Value = "TOP2"
$ws.Range("E1").Value = "TOP4"
$ws.Range("F1").Value = "RELEGATION"

# The inserted columns start out blank for every data row (2-19);
# materialize them as empty text cells.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 2).End(-4162).Row
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 3).NumberFormat = "@"
    $ws.Cells.Item($r, 4).NumberFormat = "@"
    $ws.Cells.Item($r, 5).NumberFormat = "@"
    $ws.Cells.Item($r, 6).NumberFormat = "@"
    $ws.Cells.Item($r, 3).Value = ""
    $ws.Cells.Item($r, 4).Value = ""
    $ws.Cells.Item($r, 5).Value = ""
    $ws.Cells.Item($r, 6).Value = ""
}
